# "Generate Report for Handoff"
# Two new files showed up in the handoff report:
#   7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md  (inserted before 88ec1a48...)
#   8a407dfc-124d-44cb-957a-d307955c1e31.md  (inserted after 88ec1a48..., before ff090798...)
# Every sheet's file list grows from 3 data rows to 5 data rows (rows 2-6),
# and every hyperlink + its display text is rebuilt to point at the right row.

$wb = $excel.ActiveWorkbook

$base = "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c"

$fileA = "7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md"
$fileB = "88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md"
$fileC = "8a407dfc-124d-44cb-957a-d307955c1e31.md"
$fileD = "ff090798-82ce-4771-adaf-679755eac184.md"
$fileE = ".localization-config"

$urlA = $base + "/e2e/" + $fileA
$urlB = $base + "/e2e/" + $fileB
$urlC = $base + "/e2e/" + $fileC
$urlD = $base + "/e2e/" + $fileD
$urlE = $base + "/" + $fileE

# ---- Sheet "Overview": File Name | zh-cn | de-de -------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Cells.Item(2,2).Value = "Handoff transform failed"
$ws.Cells.Item(2,3).Value = "Handoff transform failed"

$ws.Cells.Item(3,2).Value = "Handoff transform failed"
$ws.Cells.Item(3,3).Value = "Handoff transform failed"

$ws.Cells.Item(4,2).Value = "Handoff transform failed"
$ws.Cells.Item(4,3).Value = "Handoff transform failed"

$ws.Cells.Item(5,2).Value = "Handoff transform failed"
$ws.Cells.Item(5,3).Value = "Handoff transform failed"

$ws.Cells.Item(6,2).Value = "Not to be localized"
$ws.Cells.Item(6,3).Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlA, "", "", $fileA)
$ws.Hyperlinks.Add($ws.Range("A3"), $urlB, "", "", $fileB)
$ws.Hyperlinks.Add($ws.Range("A4"), $urlC, "", "", $fileC)
$ws.Hyperlinks.Add($ws.Range("A5"), $urlD, "", "", $fileD)
$ws.Hyperlinks.Add($ws.Range("A6"), $urlE, "", "", $fileE)

# ---- Sheets "zh-cn" and "de-de": identical per-language detail sheets ----
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Hyperlinks.Delete()

    $ws.Cells.Item(2,2).Value = "Handoff transform failed"
    $ws.Cells.Item(2,4).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(2,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(2,8).Value = "Ignored"

    $ws.Cells.Item(3,2).Value = "Handoff transform failed"
    $ws.Cells.Item(3,4).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(3,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(3,8).Value = "Ignored"

    $ws.Cells.Item(4,2).Value = "Handoff transform failed"
    $ws.Cells.Item(4,4).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(4,8).Value = "Ignored"

    $ws.Cells.Item(5,2).Value = "Handoff transform failed"
    $ws.Cells.Item(5,4).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(5,8).Value = "Ignored"

    $ws.Cells.Item(6,2).Value = "Not to be localized"
    $ws.Cells.Item(6,4).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(6,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(6,8).Value = "Ignored"

    $ws.Hyperlinks.Add($ws.Range("A2"), $urlA, "", "", $fileA)
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlB, "", "", $fileB)
    $ws.Hyperlinks.Add($ws.Range("A4"), $urlC, "", "", $fileC)
    $ws.Hyperlinks.Add($ws.Range("A5"), $urlD, "", "", $fileD)
    $ws.Hyperlinks.Add($ws.Range("A6"), $urlE, "", "", $fileE)

    $ws.Range("D2:D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
}
